# Add a new "2022-Q4" sheet before the current "2022-Q3" sheet, fill it
# with the quarterly fund-holding data, and update the "总计" (total)
# summary sheet with the new quarter's aggregate row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right before the existing
#    "2022-Q3" sheet (which is currently the 2nd sheet).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Populate the header row (same columns/labels as the other
#    per-quarter sheets) and copy their formatting so the new sheet
#    looks consistent with its siblings.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Copy header + index-column formatting from the "总计" sheet, which
# already uses the bold/centered/bordered style for these cells.
$total = $wb.Worksheets.Item(1)
$total.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Fill in the fund holdings data for 2022-Q4.
# ---------------------------------------------------------------------
$data = @(
    @("001852", "融通中国风1号灵活配置混合A", "23.60", "91.60", "4.60", "1.0856", 10),
    @("161606", "融通行业景气混合A/B", "17.53", "94.53", "4.62", "0.8099", 9),
    @("002989", "融通通乾研究精选灵活配置混合", "4.09", "94.13", "4.57", "0.1869", 9),
    @("011011", "融通产业趋势精选2年封闭运作混合", "2.65", "94.22", "4.54", "0.1203", 9),
    @("008382", "融通产业趋势股票", "1.88", "93.43", "5.08", "0.0955", 8),
    @("009277", "融通行业景气混合C", "1.76", "94.53", "4.62", "0.0813", 9),
    @("000916", "前海开源股息率100强等权重股票", "2.31", "91.41", "1.20", "0.0277", 7),
    @("009273", "融通中国风1号灵活配置混合C", "0.37", "91.60", "4.60", "0.0170", 10),
    @("007084", "天治转型升级混合", "0.07", "90.23", "7.95", "0.0056", 8)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 2 + $r
    $rec = $data[$r]

    # Index column (A) - numeric, styled like the rest of the sheets.
    $q4.Cells.Item($row, 1).Value = $r
    $total.Range("A2").Copy()
    $q4.Cells.Item($row, 1).PasteSpecial(-4122)

    # Fund code (B) - force text so leading zeros survive.
    $codeCell = $q4.Cells.Item($row, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $rec[0]

    # Fund name (C) - plain text.
    $q4.Cells.Item($row, 3).Value = $rec[1]

    # Fund scale / stock position / position pct / holding value (D-G)
    # - force text so the values keep their original formatting (e.g.
    # trailing zeros) exactly like the sibling sheets.
    for ($c = 0; $c -le 3; $c++) {
        $cell = $q4.Cells.Item($row, 4 + $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rec[2 + $c]
    }

    # Position rank (H) - numeric.
    $q4.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# 4. Update the "总计" (total) sheet: insert a new row for 2022-Q4 right
#    after the header, pushing the existing rows down, and renumber the
#    index column.
# ---------------------------------------------------------------------

# Make sure row 6 (the new last row) has the same formatting as the
# other index-column / data rows before shifting values into it.
$total.Range("A5:D5").Copy()
$total.Range("A6:D6").PasteSpecial(-4122)

# Shift the old quarter rows (2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4) down
# by one row.
$oldRows = $total.Range("A2:D5").Value()
$total.Range("A3:D6").Value = $oldRows

# Write the new 2022-Q4 row into the now-empty row 2.
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 9
$total.Cells.Item(2, 4).Value = 2.43

# Renumber the index column (A) 0..4.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
